# This script reproduces (on slides 2 and 3) the addition of duplicated
# picture/arrow shapes described by the target diff.
#
# PowerPoint's Shape.Left/Top/Width/Height COM properties are expressed in
# points (1 pt = 12700 EMU) and are backed by single-precision floats that
# get truncated (not rounded) when PowerPoint converts them back to EMU on
# save. Empirically, adding a small constant epsilon (0.000035 pt) before
# the conversion reliably compensates for that truncation across the whole
# range of slide coordinates, so we always land on the exact target EMU
# value. EmuToPt() below centralizes that correction.

$EMU_PER_PT = 12700.0
$EPS = 0.000035

function EmuToPt($emu) {
    return ($emu / $EMU_PER_PT) + $EPS
}

function SetBounds($shape, $x, $y, $cx, $cy) {
    $shape.LockAspectRatio = [Microsoft.Office.Core.MsoTriState]::msoFalse
    $shape.Left = EmuToPt $x
    $shape.Top = EmuToPt $y
    $shape.Width = EmuToPt $cx
    $shape.Height = EmuToPt $cy
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2: append two duplicated pictures and a duplicated right-arrow
# shape at the end of the shape tree.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$srcPic1 = $s2.Shapes("图片 1")
$newPic12 = $srcPic1.Duplicate().Item(1)
$newPic12.Name = "图片 12"
SetBounds $newPic12 414396 4234920 2298882 3130858

$srcPic2 = $s2.Shapes("图片 10")
$newPic13 = $srcPic2.Duplicate().Item(1)
$newPic13.Name = "图片 13"
SetBounds $newPic13 2757256 4243798 2276988 3123560

$srcArrow = $s2.Shapes("箭头: 右 11")
$newArrow14 = $srcArrow.Duplicate().Item(1)
$newArrow14.Name = "箭头: 右 14"
$newArrow14.Rotation = 0
SetBounds $newArrow14 2341928 5513724 848825 573250

# ---------------------------------------------------------------------
# Slide 3: reposition/rename/resize the two existing pictures and the
# existing arrow, then append a duplicated picture and a duplicated,
# rotated right-arrow shape.
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

$pic1 = $s3.Shapes("图片 1")
$pic1.Name = "图片 11"
SetBounds $pic1 3870528 244199 2284286 3108964
$pic1.Fill.Visible = [Microsoft.Office.Core.MsoTriState]::msoFalse

$pic10 = $s3.Shapes("图片 10")
$pic10.Name = "图片 9"
SetBounds $pic10 1696350 3429000 2298882 3130858

$arrow7 = $s3.Shapes("箭头: 右 7")
$arrow7.Name = "箭头: 右 13"
SetBounds $arrow7 4277219 4649786 1523217 689286

$srcPic3 = $s3.Shapes("图片 9")
$newPic12b = $srcPic3.Duplicate().Item(1)
$newPic12b.Name = "图片 12"
SetBounds $newPic12b 6030113 3429000 2276988 3123560

$srcArrow3 = $s3.Shapes("箭头: 右 13")
$newArrow15 = $srcArrow3.Duplicate().Item(1)
$newArrow15.Name = "箭头: 右 14"
$newArrow15.Rotation = 0
SetBounds $newArrow15 4499260 3900277 1026820 414143
$newArrow15.Rotation = 90
